# Natmi following Dr Hou advice
# Update Wnt5a-Mcam LR-pair results: recompute stats (now based on 3 clusters
# instead of 1) and add the two missing FAPs/sCs/ECs combinations as new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FAPs -> Wnt5a/Mcam -> ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Wnt5a"
$ws.Range("C2").Value = "Mcam"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 9.156959333333335
$ws.Range("H2").Value = 27.470878
$ws.Range("I2").Value = 0.969469463764299
$ws.Range("J2").Value = 0.9694694637642989
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 45.50099733333334
$ws.Range("N2").Value = 136.502992
$ws.Range("O2").Value = 0.5190633550775298
$ws.Range("P2").Value = 0.5190633550775298
$ws.Range("Q2").Value = 416.6507822074419
$ws.Range("R2").Value = 3749.857039866977
$ws.Range("S2").Value = 0.5032160725067107
$ws.Range("T2").Value = 0.5032160725067106

# Row 3: FAPs -> Wnt5a/Mcam -> FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Wnt5a"
$ws.Range("C3").Value = "Mcam"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 9.156959333333335
$ws.Range("H3").Value = 27.470878
$ws.Range("I3").Value = 0.969469463764299
$ws.Range("J3").Value = 0.9694694637642989
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.5123886666666667
$ws.Range("N3").Value = 1.537166
$ws.Range("O3").Value = 0.005845194523436572
$ws.Range("P3").Value = 0.005845194523436572
$ws.Range("Q3").Value = 4.691922183527557
$ws.Range("R3").Value = 42.22729965174801
$ws.Range("S3").Value = 0.005666737600234071
$ws.Range("T3").Value = 0.00566673760023407

# Row 4: FAPs -> Wnt5a/Mcam -> sCs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Wnt5a"
$ws.Range("C4").Value = "Mcam"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 9.156959333333335
$ws.Range("H4").Value = 27.470878
$ws.Range("I4").Value = 0.969469463764299
$ws.Range("J4").Value = 0.9694694637642989
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 41.64642833333333
$ws.Range("N4").Value = 124.939285
$ws.Range("O4").Value = 0.4750914503990336
$ws.Range("P4").Value = 0.4750914503990337
$ws.Range("Q4").Value = 381.3546506269145
$ws.Range("R4").Value = 3432.19185564223
$ws.Range("S4").Value = 0.4605866536573542
$ws.Range("T4").Value = 0.4605866536573542

# Row 5: sCs -> Wnt5a/Mcam -> ECs
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Wnt5a"
$ws.Range("C5").Value = "Mcam"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.288371
$ws.Range("H5").Value = 0.865113
$ws.Range("I5").Value = 0.03053053623570109
$ws.Range("J5").Value = 0.03053053623570109
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 45.50099733333334
$ws.Range("N5").Value = 136.502992
$ws.Range("O5").Value = 0.5190633550775298
$ws.Range("P5").Value = 0.5190633550775298
$ws.Range("Q5").Value = 13.12116810201067
$ws.Range("R5").Value = 118.090512918096
$ws.Range("S5").Value = 0.0158472825708191
$ws.Range("T5").Value = 0.0158472825708191

# Row 6 (new): sCs -> Wnt5a/Mcam -> FAPs
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Wnt5a"
$ws.Range("C6").Value = "Mcam"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.288371
$ws.Range("H6").Value = 0.865113
$ws.Range("I6").Value = 0.03053053623570109
$ws.Range("J6").Value = 0.03053053623570109
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.5123886666666667
$ws.Range("N6").Value = 1.537166
$ws.Range("O6").Value = 0.005845194523436572
$ws.Range("P6").Value = 0.005845194523436572
$ws.Range("Q6").Value = 0.1477580321953333
$ws.Range("R6").Value = 1.329822289758
$ws.Range("S6").Value = 0.0001784569232025018
$ws.Range("T6").Value = 0.0001784569232025018

# Row 7 (new): sCs -> Wnt5a/Mcam -> sCs
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Wnt5a"
$ws.Range("C7").Value = "Mcam"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.288371
$ws.Range("H7").Value = 0.865113
$ws.Range("I7").Value = 0.03053053623570109
$ws.Range("J7").Value = 0.03053053623570109
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 41.64642833333333
$ws.Range("N7").Value = 124.939285
$ws.Range("O7").Value = 0.4750914503990336
$ws.Range("P7").Value = 0.4750914503990337
$ws.Range("Q7").Value = 12.00962218491167
$ws.Range("R7").Value = 108.086599664205
$ws.Range("S7").Value = 0.01450479674167948
$ws.Range("T7").Value = 0.01450479674167948
